$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New donor order (name, amountUSD) for rows 2-55, matching the reorder
# performed on the underlying data (org type / flow date classification).
$data = @(
    @("United States of America, Government of", 14157836),
    @("Denmark, Government of", 699753),
    @("Austria, Government of", 1169591),
    @("Switzerland, Government of", 2791562),
    @("France, Government of", 1034819),
    @("Spain, Government of", 315948),
    @("Saudi Arabia (Kingdom of), Government of", 1690120),
    @("European Commission's Humanitarian Aid and Civil Protection Department", 2095295),
    @("Australia, Government of", 6809544),
    @("Sheikh Abdullah Al-Nouri Charity Society", 21287),
    @("Japan, Government of", 6928508),
    @("Qatar Charity", 258587),
    @("Swiss Solidarity", 8641870),
    @("Taiwan International Cooperation and Development Fund", 400000),
    @("Netherlands, Government of", 4441855),
    @("Finland, Government of", 113636),
    @("Sweden, Government of", 2897107),
    @("Liechtenstein, Government of", 177970),
    @("New Zealand, Government of", 2759429),
    @("Norway, Government of", 1231072),
    @("Central Emergency Response Fund", 14425047),
    @("Jersey Overseas Aid", 447570),
    @("US Fund for UNICEF", 560669),
    @("UNICEF National Committee/Canada", 32214),
    @("UNICEF National Committee/Japan", 156480),
    @("UNICEF National Committee/Netherlands", 258947),
    @("UNICEF National Committee/Sweden", 24750),
    @("UNICEF National Committee/Malaysia", 262506),
    @("UNICEF National Committee/Indonesia", 125267),
    @("Education Cannot Wait Fund", 1600000),
    @("UNICEF National Committee/Spain", 1557982),
    @("UNICEF National Committee/Portugal", 547021),
    @("UNICEF National Committee/Germany", 308599),
    @("UNICEF National Committee/Australia", 159820),
    @("UNICEF National Committee/Italy", 113259),
    @("UNICEF National Committee/Korea (Republic of)", 133597),
    @("UNICEF National Committee/Norway", 77872),
    @("UNICEF National Committee/United Kingdom", 119995),
    @("UNICEF National Committee/Austria", 56117),
    @("UNICEF National Committee/Belgium", 161802),
    @("UNICEF National Committee/Czech Republic", 12514),
    @("UNICEF National Committee/Denmark", 68959),
    @("UNICEF National Committee/Luxembourg", 228119),
    @("UNICEF National Committee/Switzerland", 234581),
    @("Canada, Government of", 790918),
    @("Luxembourg, Government of", 219747),
    @("Samenwerkende Hulporganisaties", 89188),
    @("United Nations Population Fund", 94142),
    @("Singapore Red Cross Society", 50000),
    @("Azerbaijan, Government of", 100000),
    @("Ireland, Government of", 501437),
    @("Monaco, Government of", 113636),
    @("Italy, Government of", 284156),
    @("Belgium, Government of", 200000)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
